# Auto-generated update for resum_diari_meteocat.xlsx
# Commit: Update automàtic: dades i banners [2026-02-21 22:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-21 22:18:27"
$ws.Range("O2").Value = "3.8 °C"
$ws.Range("E3").Value = "2026-02-21 22:18:30"
$ws.Range("H3").Value = "'41%"
$ws.Range("E4").Value = "2026-02-21 22:18:32"
$ws.Range("O4").Value = "9.2 °C"
$ws.Range("E5").Value = "2026-02-21 22:18:35"
$ws.Range("O5").Value = "3.8 °C"
$ws.Range("E6").Value = "2026-02-21 22:18:37"
$ws.Range("E7").Value = "2026-02-21 22:18:40"
$ws.Range("J7").Value = "1029.6 hPa"
$ws.Range("E8").Value = "2026-02-21 22:18:42"
$ws.Range("E9").Value = "2026-02-21 22:18:45"
$ws.Range("H9").Value = "'58%"
$ws.Range("N9").Value = "6.5 °C 21:54 TU"
$ws.Range("O9").Value = "13.0 °C"
$ws.Range("E10").Value = "2026-02-21 22:18:47"
$ws.Range("E11").Value = "2026-02-21 22:18:50"
$ws.Range("H11").Value = "'54%"
$ws.Range("O11").Value = "8.7 °C"
$ws.Range("E12").Value = "2026-02-21 22:18:52"
$ws.Range("H12").Value = "'64%"
$ws.Range("O12").Value = "12.4 °C"
$ws.Range("E13").Value = "2026-02-21 22:18:54"
$ws.Range("J13").Value = "1031.9 hPa"
$ws.Range("O13").Value = "5.2 °C"
$ws.Range("E14").Value = "2026-02-21 22:18:57"
$ws.Range("O14").Value = "11.2 °C"
$ws.Range("E15").Value = "2026-02-21 22:18:59"
$ws.Range("H15").Value = "'57%"
$ws.Range("N15").Value = "5.0 °C 21:59 TU"
$ws.Range("O15").Value = "12.8 °C"
$ws.Range("E16").Value = "2026-02-21 22:19:01"
$ws.Range("E17").Value = "2026-02-21 22:19:04"
$ws.Range("E18").Value = "2026-02-21 22:19:06"
$ws.Range("O18").Value = "8.5 °C"
$ws.Range("E19").Value = "2026-02-21 22:19:09"
$ws.Range("E20").Value = "2026-02-21 22:19:11"
$ws.Range("E21").Value = "2026-02-21 22:19:14"
$ws.Range("H21").Value = "'57%"
$ws.Range("J21").Value = "1030.8 hPa"
$ws.Range("E22").Value = "2026-02-21 22:19:16"
$ws.Range("K22").Value = "16.6 MJ/m2"
$ws.Range("E23").Value = "2026-02-21 22:19:19"
$ws.Range("O23").Value = "2.9 °C"
$ws.Range("E24").Value = "2026-02-21 22:19:21"
$ws.Range("H24").Value = "'84%"
$ws.Range("E25").Value = "2026-02-21 22:19:24"
$ws.Range("E26").Value = "2026-02-21 22:19:26"
$ws.Range("J26").Value = "1027.5 hPa"
$ws.Range("O26").Value = "9.6 °C"
$ws.Range("E27").Value = "2026-02-21 22:19:29"
$ws.Range("E28").Value = "2026-02-21 22:19:32"
$ws.Range("O28").Value = "8.1 °C"
$ws.Range("E29").Value = "2026-02-21 22:19:34"
$ws.Range("N29").Value = "5.3 °C 21:48 TU"
$ws.Range("O29").Value = "11.3 °C"
$ws.Range("E30").Value = "2026-02-21 22:19:37"
$ws.Range("O30").Value = "11.4 °C"
$ws.Range("E31").Value = "2026-02-21 22:19:39"
$ws.Range("E32").Value = "2026-02-21 22:19:42"
$ws.Range("H32").Value = "'81%"
$ws.Range("O32").Value = "4.8 °C"
$ws.Range("E33").Value = "2026-02-21 22:19:44"
$ws.Range("J33").Value = "1030.5 hPa"
$ws.Range("O33").Value = "6.6 °C"
$ws.Range("E34").Value = "2026-02-21 22:19:47"
$ws.Range("K34").Value = "14.1 MJ/m2"
$ws.Range("E35").Value = "2026-02-21 22:19:49"
$ws.Range("J35").Value = "1031.1 hPa"
$ws.Range("O35").Value = "7.5 °C"
$ws.Range("E36").Value = "2026-02-21 22:19:52"
$ws.Range("H36").Value = "'59%"
$ws.Range("N36").Value = "6.8 °C 21:44 TU"
$ws.Range("O36").Value = "13.1 °C"
$ws.Range("E37").Value = "2026-02-21 22:19:54"
$ws.Range("J37").Value = "1031.7 hPa"
$ws.Range("O37").Value = "5.6 °C"
$ws.Range("E38").Value = "2026-02-21 22:19:57"
$ws.Range("H38").Value = "'73%"
$ws.Range("E39").Value = "2026-02-21 22:19:59"
$ws.Range("E40").Value = "2026-02-21 22:20:02"
$ws.Range("H40").Value = "'54%"
$ws.Range("J40").Value = "1030.8 hPa"
$ws.Range("O40").Value = "8.4 °C"
$ws.Range("E41").Value = "2026-02-21 22:20:04"
$ws.Range("H41").Value = "'70%"
$ws.Range("E42").Value = "2026-02-21 22:20:07"
$ws.Range("H42").Value = "'76%"
$ws.Range("O42").Value = "10.6 °C"
$ws.Range("E43").Value = "2026-02-21 22:20:09"
$ws.Range("K43").Value = "15.0 MJ/m2"
$ws.Range("E44").Value = "2026-02-21 22:20:12"
$ws.Range("E45").Value = "2026-02-21 22:20:14"
$ws.Range("E46").Value = "2026-02-21 22:20:17"
$ws.Range("H46").Value = "'70%"
$ws.Range("O46").Value = "9.6 °C"
